# Update "想去人数" (interest count) figures on the 展览 (Exhibitions) sheet
# and the corresponding duplicated rows on the 全部类型 (All Types) sheet,
# matching the data refresh captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 12590   # was 12584
$ws1.Range("F6").Value  = 58      # was 57
$ws1.Range("F10").Value = 319     # was 318
$ws1.Range("F17").Value = 329     # was 328
$ws1.Range("F18").Value = 222     # was 221
$ws1.Range("F19").Value = 285     # was 284
$ws1.Range("F26").Value = 61      # was 59
$ws1.Range("F27").Value = 101     # was 100

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 12590   # was 12584
$ws4.Range("F10").Value = 58      # was 57
$ws4.Range("F14").Value = 319     # was 318
$ws4.Range("F29").Value = 329     # was 328
$ws4.Range("F31").Value = 222     # was 221
$ws4.Range("F32").Value = 285     # was 284
$ws4.Range("F43").Value = 61      # was 59
$ws4.Range("F44").Value = 101     # was 100
